$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @{
    "K8" = 1.5
    "K9" = 1.5
    "K26" = 1.5
    "K29" = 2.5
    "K32" = 1.5
    "K38" = 1.5
    "K39" = 2.5
    "K40" = 2.5
    "K44" = 1.5
    "K45" = 2.5
    "K46" = 1.5
    "K48" = 2.5
    "K55" = 3.5
    "K66" = 2.5
    "K67" = 2.5
    "K68" = 2.5
    "K70" = 1.5
    "K82" = 1.5
    "K84" = 4.5
    "K85" = 5.5
    "K89" = 0.8
    "K90" = 1.5
    "K94" = 1.5
    "K97" = 1.5
    "K98" = 2.5
    "K99" = 2.5
    "K101" = 2.5
    "K110" = 5.5
    "K125" = 2.5
    "K133" = 2.5
    "K134" = 1.5
    "K135" = 2.5
    "K137" = 1.5
    "K139" = 2.5
}

foreach ($addr in $cells.Keys) {
    $ws.Range($addr).Value = $cells[$addr]
}

$ws.Columns("K:K").Select()
